# Applies updated Leve profit/price figures per scheduled runner refresh.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 5548
$ws.Range("I2").Value = 3150.5
$ws.Range("J2").Value = 6746.75
$ws.Range("K2").Value = 3150.5
$ws.Range("L2").Value = 6746.75
$ws.Range("M2").Value = -3037.5
$ws.Range("N2").Value = -6972.75
$ws.Range("H33").Value = 310
$ws.Range("I33").Value = 265.7143
$ws.Range("K33").Value = 265.7143
$ws.Range("M33").Value = -36.71429999999998
$ws.Range("H48").Value = 1500
$ws.Range("I48").Value = 1500
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 4500
$ws.Range("L48").Value = 0
$ws.Range("M48").Value = -4208
$ws.Range("N48").ClearContents()
$ws.Range("H56").Value = 1500
$ws.Range("I56").Value = 1500
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 4500
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -3966
$ws.Range("N56").ClearContents()
$ws.Range("H86").Value = 5598.826
$ws.Range("I86").Value = 4733.2666
$ws.Range("J86").Value = 7221.75
$ws.Range("K86").Value = 4733.2666
$ws.Range("L86").Value = 7221.75
$ws.Range("M86").Value = -3610.2666
$ws.Range("N86").Value = -9467.75
$ws.Range("H89").Value = 5598.826
$ws.Range("I89").Value = 4733.2666
$ws.Range("J89").Value = 7221.75
$ws.Range("K89").Value = 23666.333
$ws.Range("L89").Value = 36108.75
$ws.Range("M89").Value = -18050.333
$ws.Range("N89").Value = -47340.75
$ws.Range("H111").Value = 991.5
$ws.Range("I111").Value = 991.5
$ws.Range("K111").Value = 2974.5
$ws.Range("M111").Value = 92.5
$ws.Range("H112").Value = 1044.4117
$ws.Range("J112").Value = 1015.9375
$ws.Range("L112").Value = 3047.8125
$ws.Range("N112").Value = -5263.8125
$ws.Range("H113").Value = 4299.5835
$ws.Range("I113").Value = 4137.125
$ws.Range("K113").Value = 4137.125
$ws.Range("M113").Value = -883.125
$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("M130").ClearContents()
$ws.Range("N130").ClearContents()
$ws.Range("H138").Value = 3984.274
$ws.Range("J138").Value = 4341.4917
$ws.Range("L138").Value = 13024.4751
$ws.Range("N138").Value = -23304.4751
$ws.Range("H141").Value = 7196.2856
$ws.Range("I141").Value = 7295
$ws.Range("K141").Value = 21885
$ws.Range("M141").Value = -16705

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1434.8334
$ws.Range("I110").Value = 1434.8334
$ws.Range("K110").Value = 1434.8334
$ws.Range("M110").Value = 610.1666
$ws.Range("H122").Value = 2207.2
$ws.Range("I122").Value = 1812.1666
$ws.Range("J122").Value = 2799.75
$ws.Range("K122").Value = 5436.4998
$ws.Range("L122").Value = 8399.25
$ws.Range("M122").Value = -2986.4998
$ws.Range("N122").Value = -13299.25

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 481.0909
$ws.Range("I5").Value = 583.8333
$ws.Range("J5").Value = 357.8
$ws.Range("K5").Value = 583.8333
$ws.Range("L5").Value = 357.8
$ws.Range("M5").Value = -470.8333
$ws.Range("N5").Value = -583.8
$ws.Range("H20").Value = 3449.5557
$ws.Range("I20").Value = 3149.5715
$ws.Range("J20").Value = 4499.5
$ws.Range("K20").Value = 3149.5715
$ws.Range("L20").Value = 4499.5
$ws.Range("M20").Value = -2902.5715
$ws.Range("N20").Value = -4993.5
$ws.Range("H74").Value = 88000
$ws.Range("J74").Value = 88000
$ws.Range("L74").Value = 88000
$ws.Range("N74").Value = -89872
$ws.Range("H77").Value = 88000
$ws.Range("J77").Value = 88000
$ws.Range("L77").Value = 264000
$ws.Range("N77").Value = -273360
$ws.Range("H94").Value = 608.3333
$ws.Range("I94").Value = 608.3333
$ws.Range("K94").Value = 608.3333
$ws.Range("M94").Value = -157.3333
$ws.Range("H107").Value = 4779.5625
$ws.Range("I107").Value = 4772.077
$ws.Range("K107").Value = 4772.077
$ws.Range("M107").Value = -2852.077

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6641.1333
$ws.Range("J31").Value = 6726.75
$ws.Range("L31").Value = 6726.75
$ws.Range("N31").Value = -7316.75
$ws.Range("H34").Value = 6641.1333
$ws.Range("J34").Value = 6726.75
$ws.Range("L34").Value = 6726.75
$ws.Range("N34").Value = -7130.75
$ws.Range("H51").Value = 38347.5
$ws.Range("I51").Value = 30090
$ws.Range("K51").Value = 30090
$ws.Range("M51").Value = -29354
$ws.Range("H59").Value = 200000
$ws.Range("I59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("M59").ClearContents()
$ws.Range("H60").Value = 4646.5
$ws.Range("I60").Value = 1575.8
$ws.Range("K60").Value = 1575.8
$ws.Range("M60").Value = -1064.8
$ws.Range("H61").Value = 38347.5
$ws.Range("I61").Value = 30090
$ws.Range("K61").Value = 30090
$ws.Range("M61").Value = -29742
$ws.Range("H100").Value = 99995
$ws.Range("J100").Value = 99995
$ws.Range("L100").Value = 99995
$ws.Range("N100").Value = -102159
$ws.Range("H107").Value = 565.875
$ws.Range("I107").Value = 517.75
$ws.Range("K107").Value = 517.75
$ws.Range("M107").Value = 1402.25
$ws.Range("H122").Value = 1709.3
$ws.Range("I122").Value = 1522.375
$ws.Range("K122").Value = 4567.125
$ws.Range("M122").Value = -2117.125
$ws.Range("H134").Value = 2786.3333
$ws.Range("I134").Value = 2668.6428
$ws.Range("J134").Value = 3198.25
$ws.Range("K134").Value = 8005.928400000001
$ws.Range("L134").Value = 9594.75
$ws.Range("M134").Value = -5470.928400000001
$ws.Range("N134").Value = -14664.75

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 16667724
$ws.Range("I4").Value = 20001242
$ws.Range("K4").Value = 60003726
$ws.Range("M4").Value = -60003614
$ws.Range("H60").Value = 473.83334
$ws.Range("I60").Value = 473.83334
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 1421.50002
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = -1170.50002
$ws.Range("N60").ClearContents()
$ws.Range("H127").Value = 1996.5
$ws.Range("J127").Value = 1996.5
$ws.Range("L127").Value = 5989.5
$ws.Range("N127").Value = -15909.5
$ws.Range("H134").Value = 13061.167
$ws.Range("I134").Value = 7126
$ws.Range("J134").Value = 18996.334
$ws.Range("K134").Value = 21378
$ws.Range("L134").Value = 56989.00199999999
$ws.Range("M134").Value = -16308
$ws.Range("N134").Value = -67129.00199999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5846
$ws.Range("J70").Value = 6015.2
$ws.Range("L70").Value = 6015.2
$ws.Range("N70").Value = -6555.2
$ws.Range("H73").Value = 5846
$ws.Range("J73").Value = 6015.2
$ws.Range("L73").Value = 6015.2
$ws.Range("N73").Value = -7887.2
$ws.Range("H80").Value = 10663
$ws.Range("J80").Value = 11940.917
$ws.Range("L80").Value = 11940.917
$ws.Range("N80").Value = -13936.917
$ws.Range("H83").Value = 10663
$ws.Range("J83").Value = 11940.917
$ws.Range("L83").Value = 59704.585
$ws.Range("N83").Value = -69688.58499999999
$ws.Range("H122").Value = 1778.4445
$ws.Range("I122").Value = 1463.25
$ws.Range("K122").Value = 4389.75
$ws.Range("M122").Value = -1939.75
$ws.Range("H132").Value = 3499.4285
$ws.Range("I132").Value = 3249.5
$ws.Range("J132").Value = 3832.6667
$ws.Range("K132").Value = 9748.5
$ws.Range("L132").Value = 11498.0001
$ws.Range("M132").Value = -7218.5
$ws.Range("N132").Value = -16558.0001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3425.3333
$ws.Range("I132").Value = 3425.3333
$ws.Range("K132").Value = 10275.9999
$ws.Range("M132").Value = -7745.999899999999
$ws.Range("H139").Value = 89650
$ws.Range("I139").Value = 89650
$ws.Range("K139").Value = 89650
$ws.Range("M139").Value = -84510

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 5556448
$ws.Range("I100").Value = 7693053.5
$ws.Range("J100").Value = 1274.8
$ws.Range("K100").Value = 15386107
$ws.Range("L100").Value = 2549.6
$ws.Range("M100").Value = -15385566
$ws.Range("N100").Value = -3631.6
$ws.Range("H107").Value = 439.14285
$ws.Range("I107").Value = 412.33334
$ws.Range("K107").Value = 1237.00002
$ws.Range("M107").Value = 682.9999800000001
$ws.Range("H126").Value = 4430.591
$ws.Range("I126").Value = 2641.9285
$ws.Range("K126").Value = 7925.7855
$ws.Range("M126").Value = -5455.7855
